$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert a new row at 3 for OptionProfilePreset (shifts everything from row 3 down by 1) ---
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "CategoryPresetSelection"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "OptionProfilePreset"
$ws.Range("D3").Value = "string"
$ws.Range("E3").Value = '"Default"'
$ws.Range("F3").Value = "Profile determines which damage types trigger DOT effects. Default/BleedOnly = physical attacks, ElementalOnly = fire/lightning/energy attacks."
$ws.Range("G3").Value = ""

# --- Step 2: Update tooltip wording for the existing damage-type multiplier rows (now rows 8-10) ---
$ws.Range("F8").Value = "DOT damage multiplier for pierce attacks. 0.0x disables DOT from pierce entirely."
$ws.Range("F9").Value = "DOT damage multiplier for slash attacks. 0.0x disables DOT from slash entirely."
$ws.Range("F10").Value = "DOT damage multiplier for blunt attacks. 0.0x disables DOT from blunt entirely."

# --- Step 3: Insert 3 new rows at 11-13 for Fire/Lightning/Energy multipliers (shifts everything from row 11 down by 3) ---
$ws.Range("A11:A13").EntireRow.Insert()

$ws.Range("A11").Value = "CategoryDamageTypeMultipliers"
$ws.Range("B11").Value = 40
$ws.Range("C11").Value = "OptionFireMultiplier"
$ws.Range("D11").Value = "float"
$ws.Range("E11").Value = "1.0f"
$ws.Range("F11").Value = "DOT damage multiplier for fire attacks. 0.0x disables DOT from fire entirely."
$ws.Range("G11").Value = ""

$ws.Range("A12").Value = "CategoryDamageTypeMultipliers"
$ws.Range("B12").Value = 50
$ws.Range("C12").Value = "OptionLightningMultiplier"
$ws.Range("D12").Value = "float"
$ws.Range("E12").Value = "1.0f"
$ws.Range("F12").Value = "DOT damage multiplier for lightning attacks. 0.0x disables DOT from lightning entirely."
$ws.Range("G12").Value = ""

$ws.Range("A13").Value = "CategoryDamageTypeMultipliers"
$ws.Range("B13").Value = 60
$ws.Range("C13").Value = "OptionEnergyMultiplier"
$ws.Range("D13").Value = "float"
$ws.Range("E13").Value = "1.0f"
$ws.Range("F13").Value = "DOT damage multiplier for energy attacks. 0.0x disables DOT from energy entirely."
$ws.Range("G13").Value = ""
